# Update "想去人数" (want-to-go count) values in column F
# for both the "展览" and "全部类型" worksheets, which carry
# duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

# Row -> New Value map for column F
$updates = @{
    7  = 1330
    8  = 1566
    10 = 437
    12 = 176
    13 = 166
    14 = 73
    15 = 117
    16 = 281
    17 = 321
    18 = 332
    19 = 1767
    22 = 181
    23 = 689
    25 = 346
    26 = 4262
    29 = 1120
    32 = 625
    34 = 317
    35 = 50
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
